$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Columns("N:N").Insert()
$wsRepay.Columns("N:N").ColumnWidth = $wsRepay.Columns("M:M").ColumnWidth

$wsRepay.Activate()
$wsRepay.Range("R7").Select()
